$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9112548232078552
$ws.Range("B1").Value = 1.142139673233032
$ws.Range("C1").Value = 1.653894901275635
$ws.Range("D1").Value = 4.647810935974121
$ws.Range("E1").Value = 2.620930433273315
